# "scale data push and pull"
#
# The "value" column (D) was recorded in units of 万元 (ten-thousand yuan).
# Rescale every numeric entry in column D (rows 2-33) by a factor of 10000
# so the figures are pushed/pulled in plain yuan -- i.e. shift the decimal
# point four places to the right. Cells with no value (e.g. row 29) are
# left untouched.
#
# The rescale is done as an exact decimal-string shift (rather than a
# floating point multiplication) so the result is the same scaled decimal
# value the source data pipeline produced, not a value perturbed by an
# extra floating-point rounding step introduced by computing
# old_value * 10000 in binary floating point.

function Shift-DecimalPoint {
    param(
        [string]$numberText,
        [int]$places
    )

    if ($numberText -match '^(-?)(\d*)\.?(\d*)$') {
        $sign = $matches[1]
        $intPart = $matches[2]
        $fracPart = $matches[3]
    } else {
        throw "Shift-DecimalPoint: not a plain decimal number: $numberText"
    }

    if ($intPart -eq '') { $intPart = '0' }

    # Make sure there are enough fractional digits to shift out of.
    while ($fracPart.Length -lt $places) {
        $fracPart = "{0}{1}" -f $fracPart, '0'
    }

    $movedDigits = $fracPart.Substring(0, $places)
    $remainingFrac = $fracPart.Substring($places)

    $newInt = "{0}{1}" -f $intPart, $movedDigits
    $newInt = $newInt.TrimStart('0')
    if ($newInt -eq '') { $newInt = '0' }

    $result = "{0}{1}" -f $sign, $newInt
    if ($remainingFrac -ne '') {
        $result = "{0}.{1}" -f $result, $remainingFrac
    }

    return $result
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 33
$valueCol = 4   # column D

for ($row = $firstRow; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, $valueCol)
    $raw = $cell.Formula

    if ($raw -ne $null -and $raw -ne "") {
        $cell.Value2 = Shift-DecimalPoint $raw 4
    }
}
